$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $found = $ws.UsedRange.Find("congenital")
    if ($found -ne $null) {
        $found.Value = "misc_long_term"
    }
}
